$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1472596.2
$ws.Range("J17").Value = 1517199.8
$ws.Range("L17").Value = 4551599.4
$ws.Range("N17").Value = -4551935.4
$ws.Range("H33").Value = 95.166664
$ws.Range("I33").Value = 94.72727
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 94.72727
$ws.Range("L33").Value = 100
$ws.Range("M33").Value = 134.27273
$ws.Range("N33").Value = -558
$ws.Range("H51").Value = 2600
$ws.Range("I51").Value = 2833.3333
$ws.Range("K51").Value = 2833.3333
$ws.Range("M51").Value = -2349.3333
$ws.Range("H62").Value = 2693.4348
$ws.Range("I62").Value = 2228.7856
$ws.Range("K62").Value = 2228.7856
$ws.Range("M62").Value = -1604.7856
$ws.Range("H64").Value = 3612.5
$ws.Range("I64").Value = 2900
$ws.Range("J64").Value = 3850
$ws.Range("K64").Value = 2900
$ws.Range("L64").Value = 3850
$ws.Range("M64").Value = -2652
$ws.Range("N64").Value = -4346
$ws.Range("H65").Value = 2693.4348
$ws.Range("I65").Value = 2228.7856
$ws.Range("K65").Value = 11143.928
$ws.Range("M65").Value = -8023.928
$ws.Range("H67").Value = 3612.5
$ws.Range("I67").Value = 2900
$ws.Range("J67").Value = 3850
$ws.Range("K67").Value = 2900
$ws.Range("L67").Value = 3850
$ws.Range("M67").Value = -2042
$ws.Range("N67").Value = -5566
$ws.Range("H76").Value = 3400
$ws.Range("I76").Value = 3400
$ws.Range("K76").Value = 3400
$ws.Range("M76").Value = -3085
$ws.Range("H79").Value = 3400
$ws.Range("I79").Value = 3400
$ws.Range("K79").Value = 3400
$ws.Range("M79").Value = -2308
$ws.Range("H106").Value = 1880.3226
$ws.Range("I106").Value = 1254.0869
$ws.Range("K106").Value = 1254.0869
$ws.Range("M106").Value = -623.0869
$ws.Range("H116").Value = 3443.25
$ws.Range("I116").Value = 1359.4546
$ws.Range("K116").Value = 1359.4546
$ws.Range("M116").Value = 2082.5454
$ws.Range("H137").Value = 2580.44
$ws.Range("I137").Value = 2462.4285
$ws.Range("K137").Value = 7387.2855
$ws.Range("M137").Value = -4837.2855
$ws.Range("H138").Value = 1567.725
$ws.Range("J138").Value = 2290.5334
$ws.Range("L138").Value = 6871.600199999999
$ws.Range("N138").Value = -17151.6002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4065.918
$ws.Range("I32").Value = 4105.3857
$ws.Range("K32").Value = 4105.3857
$ws.Range("M32").Value = -3818.3857
$ws.Range("H61").Value = 2097.923
$ws.Range("I61").Value = 1772.75
$ws.Range("K61").Value = 1772.75
$ws.Range("M61").Value = -1560.75
$ws.Range("H74").Value = 90910050
$ws.Range("I74").Value = 100000856
$ws.Range("K74").Value = 100000856
$ws.Range("M74").Value = -99999982
$ws.Range("H77").Value = 90910050
$ws.Range("I77").Value = 100000856
$ws.Range("K77").Value = 500004280
$ws.Range("M77").Value = -499999912
$ws.Range("H132").Value = 11790.633
$ws.Range("I132").Value = 1403.0714
$ws.Range("K132").Value = 4209.2142
$ws.Range("M132").Value = -1679.2142
$ws.Range("H136").Value = 2097.923
$ws.Range("I136").Value = 1772.75
$ws.Range("K136").Value = 5318.25
$ws.Range("M136").Value = -2768.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5161.864
$ws.Range("I134").Value = 5518.05
$ws.Range("K134").Value = 16554.15
$ws.Range("M134").Value = -14019.15

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13281.897
$ws.Range("I31").Value = 22610.525
$ws.Range("K31").Value = 22610.525
$ws.Range("M31").Value = -22315.525
$ws.Range("H34").Value = 13281.897
$ws.Range("I34").Value = 22610.525
$ws.Range("K34").Value = 22610.525
$ws.Range("M34").Value = -22408.525
$ws.Range("H58").Value = 11862.761
$ws.Range("I58").Value = 1078.4166
$ws.Range("J58").Value = 23627.5
$ws.Range("K58").Value = 1078.4166
$ws.Range("L58").Value = 23627.5
$ws.Range("M58").Value = -875.4166
$ws.Range("N58").Value = -24033.5
$ws.Range("H62").Value = 125003870
$ws.Range("I62").Value = 200004400
$ws.Range("J62").Value = 2999.6667
$ws.Range("K62").Value = 200004400
$ws.Range("L62").Value = 2999.6667
$ws.Range("M62").Value = -200003776
$ws.Range("N62").Value = -4247.6667
$ws.Range("H65").Value = 125003870
$ws.Range("I65").Value = 200004400
$ws.Range("J65").Value = 2999.6667
$ws.Range("K65").Value = 1000022000
$ws.Range("L65").Value = 14998.3335
$ws.Range("M65").Value = -1000018880
$ws.Range("N65").Value = -21238.3335
$ws.Range("H86").Value = 7946142
$ws.Range("I86").Value = 2973
$ws.Range("K86").Value = 2973
$ws.Range("M86").Value = -1850
$ws.Range("H89").Value = 7946142
$ws.Range("I89").Value = 2973
$ws.Range("K89").Value = 14865
$ws.Range("M89").Value = -9249
$ws.Range("H107").Value = 543.46155
$ws.Range("I107").Value = 367.36
$ws.Range("J107").Value = 857.9286
$ws.Range("K107").Value = 367.36
$ws.Range("L107").Value = 857.9286
$ws.Range("M107").Value = 1552.64
$ws.Range("N107").Value = -4697.9286
$ws.Range("H134").Value = 1045.0758
$ws.Range("I134").Value = 768.119
$ws.Range("K134").Value = 2304.357
$ws.Range("M134").Value = 230.643
$ws.Range("H136").Value = 11862.761
$ws.Range("I136").Value = 1078.4166
$ws.Range("J136").Value = 23627.5
$ws.Range("K136").Value = 3235.2498
$ws.Range("L136").Value = 70882.5
$ws.Range("M136").Value = -685.2498000000001
$ws.Range("N136").Value = -75982.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 500.72
$ws.Range("I113").Value = 423
$ws.Range("K113").Value = 1269
$ws.Range("M113").Value = 901
$ws.Range("H131").Value = 789.51
$ws.Range("I131").Value = 333
$ws.Range("J131").Value = 794.1212
$ws.Range("K131").Value = 999
$ws.Range("L131").Value = 2382.3636
$ws.Range("M131").Value = 4041
$ws.Range("N131").Value = -12462.3636

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6139.154
$ws.Range("I70").Value = 6542.857
$ws.Range("K70").Value = 6542.857
$ws.Range("M70").Value = -6272.857
$ws.Range("H73").Value = 6139.154
$ws.Range("I73").Value = 6542.857
$ws.Range("K73").Value = 6542.857
$ws.Range("M73").Value = -5606.857
$ws.Range("H80").Value = 3399.348
$ws.Range("I80").Value = 3120.5
$ws.Range("J80").Value = 3613.8462
$ws.Range("K80").Value = 3120.5
$ws.Range("L80").Value = 3613.8462
$ws.Range("M80").Value = -2122.5
$ws.Range("N80").Value = -5609.8462
$ws.Range("H83").Value = 3399.348
$ws.Range("I83").Value = 3120.5
$ws.Range("J83").Value = 3613.8462
$ws.Range("K83").Value = 15602.5
$ws.Range("L83").Value = 18069.231
$ws.Range("M83").Value = -10610.5
$ws.Range("N83").Value = -28053.231
$ws.Range("H113").Value = 2991.3333
$ws.Range("I113").Value = 2306.3635
$ws.Range("J113").Value = 4875
$ws.Range("K113").Value = 2306.3635
$ws.Range("L113").Value = 4875
$ws.Range("M113").Value = -136.3634999999999
$ws.Range("N113").Value = -9215
$ws.Range("H132").Value = 17076.861
$ws.Range("I132").Value = 3068.276
$ws.Range("K132").Value = 9204.828
$ws.Range("M132").Value = -6674.828

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4167.294
$ws.Range("I40").Value = 2678.4285
$ws.Range("J40").Value = 5209.5
$ws.Range("K40").Value = 2678.4285
$ws.Range("L40").Value = 5209.5
$ws.Range("M40").Value = -2542.4285
$ws.Range("N40").Value = -5481.5
$ws.Range("H132").Value = 2089.9768
$ws.Range("I132").Value = 941.8214
$ws.Range("J132").Value = 4233.2
$ws.Range("K132").Value = 2825.4642
$ws.Range("L132").Value = 12699.6
$ws.Range("M132").Value = -295.4642000000003
$ws.Range("N132").Value = -17759.6
$ws.Range("H136").Value = 24806.38
$ws.Range("I136").Value = 30408.53
$ws.Range("J136").Value = 997.25
$ws.Range("K136").Value = 91225.59
$ws.Range("L136").Value = 2991.75
$ws.Range("M136").Value = -88675.59
$ws.Range("N136").Value = -8091.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1516.6666
$ws.Range("I96").Value = 1275
$ws.Range("K96").Value = 1275
$ws.Range("M96").Value = 98
$ws.Range("H132").Value = 1323.1724
$ws.Range("I132").Value = 1021.5455
$ws.Range("K132").Value = 3064.6365
$ws.Range("M132").Value = -534.6364999999996
$ws.Range("H136").Value = 31251940
$ws.Range("I136").Value = 47620760
$ws.Range("J136").Value = 2373.0908
$ws.Range("K136").Value = 142862280
$ws.Range("L136").Value = 7119.2724
$ws.Range("M136").Value = -142859730
$ws.Range("N136").Value = -12219.2724
